# Auto-generated edit script: update cryptos list values (prices & 1h volume %)
# Rows 43/44 and 46/47/48 were reordered (coin rows swapped) per upstream source.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.843.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.632.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.622.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("E10").Value = "  -4.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.620"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.28%  "

$ws.Range("E14").Value = "  -1.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.218.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "668.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.60%  "

$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.628.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.837.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.09%  "

$ws.Range("E20").Value = "  -0.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.943"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.96%  "

$ws.Range("E24").Value = "  -3.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.20%  "

$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "574.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.108"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("E41").Value = "  -0.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.549.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.38%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.140"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.94%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.345"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.49%  "

$ws.Range("E45").Value = "  -4.06%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.15%  "

$ws.Range("B47").Value = "PEPE"
$ws.Range("C47").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0730"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.21%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.135"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.30%  "
